# Swap the order of "System" and the email address in column G
# "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
# "System, admin@admin.com"    -> "admin@admin.com, System"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G
    $val = $cell.Value()
    if ($val -eq "System, dnasr281@gmail.com") {
        $cell.Value = "dnasr281@gmail.com, System"
    } elseif ($val -eq "System, admin@admin.com") {
        $cell.Value = "admin@admin.com, System"
    }
}
